# aggiornamento fino a 1/09/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44432, 0, 12, 66.74824785849371),
    @(44433, 2, 14, 77.87295583490933),
    @(44434, 3, 10, 55.6235398820781),
    @(44435, 0, 6, 33.37412392924686),
    @(44436, 4, 10, 55.6235398820781),
    @(44437, 8, 18, 100.1223717877406),
    @(44438, 5, 22, 122.3717877405718),
    @(44439, 3, 25, 139.0588497051953),
    @(44440, 0, 23, 127.9341417287796)
)

$startRow = 358
$formatSourceRow = $startRow - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Mirror the date-cell formatting (style index) from the last existing row
    $ws.Cells.Item($formatSourceRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
